# Weekly update: insert a new week's price record for
# "Terminal Hortofrutícola Agro Chillán" - Ajo (Chino) above the current
# row 128, shifting the existing rows 128-143 down to 129-144.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 128 (shifts rows 128:143 down to 129:144)
$ws.Rows(128).Insert()

# Populate the newly inserted row 128 with this week's record.
$ws.Cells.Item(128, 1).Value = 7
$ws.Cells.Item(128, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(128, 3).Value = "Ñuble"
$ws.Cells.Item(128, 4).Value = 44505
$ws.Cells.Item(128, 5).Value = 16
$ws.Cells.Item(128, 6).Value = 100112003
$ws.Cells.Item(128, 7).Value = "Ajo"
$ws.Cells.Item(128, 8).Value = "Chino"
$ws.Cells.Item(128, 9).Value = "Primera"
$ws.Cells.Item(128, 10).Value = 100
$ws.Cells.Item(128, 11).Value = 16000
$ws.Cells.Item(128, 12).Value = 17000
$ws.Cells.Item(128, 13).Value = 16500
$ws.Cells.Item(128, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(128, 15).Value = "China"
$ws.Cells.Item(128, 16).Value = 1650
$ws.Cells.Item(128, 17).Value = 10
$ws.Cells.Item(128, 18).Value = "Hortaliza"
